$d = $word.ActiveDocument

$d.Content.Find.Execute("409÷7=58, 3", $true, $false, $false, $false, $false, $true, 1, $false, "994÷2=497, 0", 2) | Out-Null
$d.Content.Find.Execute("642÷9=71, 3", $true, $false, $false, $false, $false, $true, 1, $false, "242÷5=48, 2", 2) | Out-Null
$d.Content.Find.Execute("465÷4=116, 1", $true, $false, $false, $false, $false, $true, 1, $false, "645÷8=80, 5", 2) | Out-Null
$d.Content.Find.Execute("680÷5=136, 0", $true, $false, $false, $false, $false, $true, 1, $false, "767÷8=95, 7", 2) | Out-Null
$d.Content.Find.Execute("269÷4=67, 1", $true, $false, $false, $false, $false, $true, 1, $false, "679÷9=75, 4", 2) | Out-Null
$d.Content.Find.Execute("955÷9=106, 1", $true, $false, $false, $false, $false, $true, 1, $false, "555÷7=79, 2", 2) | Out-Null
$d.Content.Find.Execute("970÷2=485, 0", $true, $false, $false, $false, $false, $true, 1, $false, "432÷5=86, 2", 2) | Out-Null
$d.Content.Find.Execute("335÷9=37, 2", $true, $false, $false, $false, $false, $true, 1, $false, "611÷2=305, 1", 2) | Out-Null
$d.Content.Find.Execute("127÷8=15, 7", $true, $false, $false, $false, $false, $true, 1, $false, "240÷8=30, 0", 2) | Out-Null
$d.Content.Find.Execute("279÷9=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "956÷7=136, 4", 2) | Out-Null
$d.Content.Find.Execute("529÷2=264, 1", $true, $false, $false, $false, $false, $true, 1, $false, "237÷9=26, 3", 2) | Out-Null
$d.Content.Find.Execute("103÷9=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "356÷8=44, 4", 2) | Out-Null
$d.Content.Find.Execute("769÷7=109, 6", $true, $false, $false, $false, $false, $true, 1, $false, "298÷4=74, 2", 2) | Out-Null
$d.Content.Find.Execute("703÷8=87, 7", $true, $false, $false, $false, $false, $true, 1, $false, "781÷8=97, 5", 2) | Out-Null
$d.Content.Find.Execute("681÷9=75, 6", $true, $false, $false, $false, $false, $true, 1, $false, "123÷4=30, 3", 2) | Out-Null
$d.Content.Find.Execute("251÷5=50, 1", $true, $false, $false, $false, $false, $true, 1, $false, "267÷8=33, 3", 2) | Out-Null
$d.Content.Find.Execute("861÷2=430, 1", $true, $false, $false, $false, $false, $true, 1, $false, "926÷8=115, 6", 2) | Out-Null
$d.Content.Find.Execute("274÷4=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "365÷3=121, 2", 2) | Out-Null
$d.Content.Find.Execute("905÷5=181, 0", $true, $false, $false, $false, $false, $true, 1, $false, "150÷2=75, 0", 2) | Out-Null
$d.Content.Find.Execute("684÷3=228, 0", $true, $false, $false, $false, $false, $true, 1, $false, "255÷5=51, 0", 2) | Out-Null
$d.Content.Find.Execute("965÷6=160, 5", $true, $false, $false, $false, $false, $true, 1, $false, "559÷9=62, 1", 2) | Out-Null
$d.Content.Find.Execute("348÷4=87, 0", $true, $false, $false, $false, $false, $true, 1, $false, "248÷5=49, 3", 2) | Out-Null
$d.Content.Find.Execute("174÷9=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "884÷6=147, 2", 2) | Out-Null
$d.Content.Find.Execute("527÷3=175, 2", $true, $false, $false, $false, $false, $true, 1, $false, "730÷3=243, 1", 2) | Out-Null
$d.Content.Find.Execute("196÷2=98, 0", $true, $false, $false, $false, $false, $true, 1, $false, "558÷6=93, 0", 2) | Out-Null
